# Weekly update: push two new price records to the top of the table
# (pushing the rest of the data rows down by 2), and backfill the two
# new rows with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 5:96 down to 7:98, inserting two blank rows
# at 5:6 (copies formatting, e.g. the date style on column D, from the
# row above - matching native Excel "insert row" behaviour).
$ws.Range("A5:R6").Insert()

# New row 5
$ws.Cells.Item(5, 1).Value = 10
$ws.Cells.Item(5, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(5, 3).Value = "La Araucanía"
$ws.Cells.Item(5, 4).Value = 44756
$ws.Cells.Item(5, 5).Value = 9
$ws.Cells.Item(5, 6).Value = 100112035
$ws.Cells.Item(5, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(5, 11).Value = 25000
$ws.Cells.Item(5, 12).Value = 25000
$ws.Cells.Item(5, 13).Value = 25000
$ws.Cells.Item(5, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(5, 15).Value = "Calera"
$ws.Cells.Item(5, 16).Value = 2500
$ws.Cells.Item(5, 17).Value = 10
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# New row 6
$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(6, 3).Value = "La Araucanía"
$ws.Cells.Item(6, 4).Value = 44756
$ws.Cells.Item(6, 5).Value = 9
$ws.Cells.Item(6, 6).Value = 100112035
$ws.Cells.Item(6, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 26000
$ws.Cells.Item(6, 12).Value = 26000
$ws.Cells.Item(6, 13).Value = 26000
$ws.Cells.Item(6, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(6, 15).Value = "Región Metropolitana"
$ws.Cells.Item(6, 16).Value = 2600
$ws.Cells.Item(6, 17).Value = 10
$ws.Cells.Item(6, 18).Value = "Hortaliza"
